$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- June Bank Statement (rows 11-13) ---
# Add a Debit Amount on the statement summary row and turn the
# "Remaining Balance" cell into a real formula instead of a hard value.
$ws.Range("K13").Value = -3145
$ws.Range("L13").Formula = "=SUM(I13:K13)"

# --- June Expense Report detail rows (14-17) ---
$ws.Range("D14").Value = 10
$ws.Range("A15").Value = 2705
$ws.Range("D15").Value = 10
$ws.Range("A16").Value = 320
$ws.Range("D16").Value = 20

# Totals row now also sums the Travel column
$ws.Range("D18").Formula = "=SUM(D16,D15,D14)"

# Recalculate so the cached formula results (A18, F18, L13, D18) are fresh
$wb.Application.Calculate()

# Restore the last-used selection
$ws.Range("J14").Select()
